# Commit: "using arrival time instead of admission in the text, put all plots in the paper"
#
# Effective spreadsheet changes in table.xlsx / Hourly sheet:
#  - Row 2 (a spacer row that only carried the second half of a
#    two-line "Multiple / seasonality" column header) is removed.
#    All the data rows below it shift up by one.
#  - The now-orphaned "seasonality" text is folded back into the I1
#    header so the column still reads " Multiple seasonality".
#  - Column I (which used to be an overflow-only default column) gets
#    an explicit width so the longer header text fits.
#  - Row 1's height is doubled to accommodate the longer wrapped text.
#  - The stray empty cell that ends up at B14 after the shift is
#    cleared out entirely (it carried no value to begin with).
#  - The selected cell/filter bookkeeping is refreshed to match the
#    new, smaller table extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the spacer row (old row 2). Everything below moves up by one.
$ws.Rows.Item(2).Delete()

# The text that used to live in (old) I2 as "seasonality" now belongs
# in the I1 header, right after "Multiple".
$ws.Range("I1").Value = " Multiple seasonality                                                                           "

# Row 1 needs to be taller to fit the combined / wrapped header text.
$ws.Cells.Item(1, 1).EntireRow.RowHeight = 88

# Column I needs an explicit width to accommodate the new header text.
# (ColumnWidth is specified in characters; Excel stores a few extra
# points of internal padding, hence the offset below nets out to a
# stored column width of 39.)
$ws.Columns.Item(9).ColumnWidth = 38.16666666666667

# The empty placeholder cell that shifts up into B14 had no content or
# special formatting to preserve, so it is cleared completely.
$ws.Cells.Item(14, 2).Clear()

# Keep the hidden filter-database bookmark in sync with the new,
# one-row-shorter table extent.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Hourly!`$A`$1:`$J`$14"

# Reflect where the user's selection ended up after the edit.
$ws.Range("I1").Select() | Out-Null
